$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")
$ws.Name = "Schema Overview"
